$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (pushes RAlt..AppsKey etc. down by one row)
$ws.Rows.Item(6).Insert()

# Fill in the new PrintScreen row (row 6) with the same formatting used by similar rows
$ws.Range("A6").Value = "PrintScreen"
$ws.Range("B6").Value = 55
$ws.Range("C6").Value = 311
$ws.Range("D6").Value = 311
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 3

# Match the style used for analogous rows (same fill/border as NumpadDiv/RAlt/Pause/Ctrl+Break/Lwin/Rwin/AppsKey)
$ws.Range("A6:F6").Style = $ws.Range("A4:F4").Style

# Update the active selection/cell to A7, matching the saved view state
$ws.Range("A7").Select()
